# "Duck shoot changes to electronics design"
# The Motor RPM input (C2) changes from 130 to 270 rpm. The two dependent
# formulas (C7 = C2/60 "Num Rotation per S" and C8 = C6/C7 "Est Time Taken")
# recalculate automatically from that single input edit. The active
# selection also moves onto the edited cell, C2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 270

$ws.Range("C2").Select()
